$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 0.2169893333333333
$ws.Range("H2").Value2 = 0.650968
$ws.Range("I2").Value2 = 0.2138091362408864
$ws.Range("J2").Value2 = 0.2138091362408864
$ws.Range("M2").Value2 = 57.35848733333334
$ws.Range("N2").Value2 = 172.075462
$ws.Range("O2").Value2 = 0.261658309594631
$ws.Range("P2").Value2 = 0.261658309594631
$ws.Range("Q2").Value2 = 12.44617992746845
$ws.Range("R2").Value2 = 112.015619347216
$ws.Range("S2").Value2 = 0.05594493716467849
$ws.Range("T2").Value2 = 0.05594493716467849
$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 0.2169893333333333
$ws.Range("H3").Value2 = 0.650968
$ws.Range("I3").Value2 = 0.2138091362408864
$ws.Range("J3").Value2 = 0.2138091362408864
$ws.Range("O3").Value2 = 0.2957894889638607
$ws.Range("P3").Value2 = 0.2957894889638607
$ws.Range("Q3").Value2 = 14.06968196806578
$ws.Range("R3").Value2 = 126.627137712592
$ws.Range("S3").Value2 = 0.06324249514449624
$ws.Range("T3").Value2 = 0.06324249514449624
$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 0.2169893333333333
$ws.Range("H4").Value2 = 0.650968
$ws.Range("I4").Value2 = 0.2138091362408864
$ws.Range("J4").Value2 = 0.2138091362408864
$ws.Range("M4").Value2 = 29.294891
$ws.Range("N4").Value2 = 87.88467299999999
$ws.Range("O4").Value2 = 0.1336376186888105
$ws.Range("P4").Value2 = 0.1336376186888105
$ws.Range("Q4").Value2 = 6.356678868162666
$ws.Range("R4").Value2 = 57.21010981346399
$ws.Range("S4").Value2 = 0.02857294382114351
$ws.Range("T4").Value2 = 0.02857294382114351
$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 0.2169893333333333
$ws.Range("H5").Value2 = 0.650968
$ws.Range("I5").Value2 = 0.2138091362408864
$ws.Range("J5").Value2 = 0.2138091362408864
$ws.Range("M5").Value2 = 67.71760166666667
$ws.Range("N5").Value2 = 203.152805
$ws.Range("O5").Value2 = 0.3089145827526977
$ws.Range("P5").Value2 = 0.3089145827526977
$ws.Range("Q5").Value2 = 14.69399724058222
$ws.Range("R5").Value2 = 132.24597516524
$ws.Range("S5").Value2 = 0.06604876011056812
$ws.Range("T5").Value2 = 0.06604876011056812
$ws.Range("I6").Value2 = 0.3272432505578689
$ws.Range("J6").Value2 = 0.3272432505578689
$ws.Range("M6").Value2 = 57.35848733333334
$ws.Range("N6").Value2 = 172.075462
$ws.Range("O6").Value2 = 0.261658309594631
$ws.Range("P6").Value2 = 0.261658309594631
$ws.Range("Q6").Value2 = 19.04936546726489
$ws.Range("R6").Value2 = 171.444289205384
$ws.Range("S6").Value2 = 0.08562591576722427
$ws.Range("T6").Value2 = 0.08562591576722427
$ws.Range("I7").Value2 = 0.3272432505578689
$ws.Range("J7").Value2 = 0.3272432505578689
$ws.Range("O7").Value2 = 0.2957894889638607
$ws.Range("P7").Value2 = 0.2957894889638607
$ws.Range("S7").Value2 = 0.09679511384938466
$ws.Range("T7").Value2 = 0.09679511384938466
$ws.Range("I8").Value2 = 0.3272432505578689
$ws.Range("J8").Value2 = 0.3272432505578689
$ws.Range("M8").Value2 = 29.294891
$ws.Range("N8").Value2 = 87.88467299999999
$ws.Range("O8").Value2 = 0.1336376186888105
$ws.Range("P8").Value2 = 0.1336376186888105
$ws.Range("Q8").Value2 = 9.729145779937332
$ws.Range("R8").Value2 = 87.56231201943599
$ws.Range("S8").Value2 = 0.04373200873653936
$ws.Range("T8").Value2 = 0.04373200873653936
$ws.Range("I9").Value2 = 0.3272432505578689
$ws.Range("J9").Value2 = 0.3272432505578689
$ws.Range("M9").Value2 = 67.71760166666667
$ws.Range("N9").Value2 = 203.152805
$ws.Range("O9").Value2 = 0.3089145827526977
$ws.Range("P9").Value2 = 0.3089145827526977
$ws.Range("Q9").Value2 = 22.48973783458445
$ws.Range("R9").Value2 = 202.40764051126
$ws.Range("S9").Value2 = 0.1010902122047206
$ws.Range("T9").Value2 = 0.1010902122047206
$ws.Range("G10").Value2 = 0.4136746666666666
$ws.Range("H10").Value2 = 1.241024
$ws.Range("I10").Value2 = 0.4076118480389355
$ws.Range("J10").Value2 = 0.4076118480389355
$ws.Range("M10").Value2 = 57.35848733333334
$ws.Range("N10").Value2 = 172.075462
$ws.Range("O10").Value2 = 0.261658309594631
$ws.Range("P10").Value2 = 0.261658309594631
$ws.Range("Q10").Value2 = 23.72775312812089
$ws.Range("R10").Value2 = 213.549778153088
$ws.Range("S10").Value2 = 0.1066550271286115
$ws.Range("T10").Value2 = 0.1066550271286115
$ws.Range("G11").Value2 = 0.4136746666666666
$ws.Range("H11").Value2 = 1.241024
$ws.Range("I11").Value2 = 0.4076118480389355
$ws.Range("J11").Value2 = 0.4076118480389355
$ws.Range("O11").Value2 = 0.2957894889638607
$ws.Range("P11").Value2 = 0.2957894889638607
$ws.Range("Q11").Value2 = 26.82284381833956
$ws.Range("R11").Value2 = 241.405594365056
$ws.Range("S11").Value2 = 0.1205673002270516
$ws.Range("T11").Value2 = 0.1205673002270516
$ws.Range("G12").Value2 = 0.4136746666666666
$ws.Range("H12").Value2 = 1.241024
$ws.Range("I12").Value2 = 0.4076118480389355
$ws.Range("J12").Value2 = 0.4076118480389355
$ws.Range("M12").Value2 = 29.294891
$ws.Range("N12").Value2 = 87.88467299999999
$ws.Range("O12").Value2 = 0.1336376186888105
$ws.Range("P12").Value2 = 0.1336376186888105
$ws.Range("Q12").Value2 = 12.11855426946133
$ws.Range("R12").Value2 = 109.066988425152
$ws.Range("S12").Value2 = 0.05447227672126863
$ws.Range("T12").Value2 = 0.05447227672126864
$ws.Range("G13").Value2 = 0.4136746666666666
$ws.Range("H13").Value2 = 1.241024
$ws.Range("I13").Value2 = 0.4076118480389355
$ws.Range("J13").Value2 = 0.4076118480389355
$ws.Range("M13").Value2 = 67.71760166666667
$ws.Range("N13").Value2 = 203.152805
$ws.Range("O13").Value2 = 0.3089145827526977
$ws.Range("P13").Value2 = 0.3089145827526977
$ws.Range("Q13").Value2 = 28.01305629692444
$ws.Range("R13").Value2 = 252.11750667232
$ws.Range("S13").Value2 = 0.1259172439620038
$ws.Range("T13").Value2 = 0.1259172439620038
$ws.Range("G14").Value2 = 0.05209933333333334
$ws.Range("H14").Value2 = 0.156298
$ws.Range("I14").Value2 = 0.05133576516230915
$ws.Range("J14").Value2 = 0.05133576516230916
$ws.Range("M14").Value2 = 57.35848733333334
$ws.Range("N14").Value2 = 172.075462
$ws.Range("O14").Value2 = 0.261658309594631
$ws.Range("P14").Value2 = 0.261658309594631
$ws.Range("Q14").Value2 = 2.988338951075112
$ws.Range("R14").Value2 = 26.895050559676
$ws.Range("S14").Value2 = 0.01343242953411676
$ws.Range("T14").Value2 = 0.01343242953411676
$ws.Range("G15").Value2 = 0.05209933333333334
$ws.Range("H15").Value2 = 0.156298
$ws.Range("I15").Value2 = 0.05133576516230915
$ws.Range("J15").Value2 = 0.05133576516230916
$ws.Range("O15").Value2 = 0.2957894889638607
$ws.Range("P15").Value2 = 0.2957894889638607
$ws.Range("Q15").Value2 = 3.378143245512445
$ws.Range("R15").Value2 = 30.40328920961201
$ws.Range("S15").Value2 = 0.01518457974292819
$ws.Range("T15").Value2 = 0.01518457974292819
$ws.Range("G16").Value2 = 0.05209933333333334
$ws.Range("H16").Value2 = 0.156298
$ws.Range("I16").Value2 = 0.05133576516230915
$ws.Range("J16").Value2 = 0.05133576516230916
$ws.Range("M16").Value2 = 29.294891
$ws.Range("N16").Value2 = 87.88467299999999
$ws.Range("O16").Value2 = 0.1336376186888105
$ws.Range("P16").Value2 = 0.1336376186888105
$ws.Range("Q16").Value2 = 1.526244291172667
$ws.Range("R16").Value2 = 13.736198620554
$ws.Range("S16").Value2 = 0.006860389409858992
$ws.Range("T16").Value2 = 0.006860389409858993
$ws.Range("G17").Value2 = 0.05209933333333334
$ws.Range("H17").Value2 = 0.156298
$ws.Range("I17").Value2 = 0.05133576516230915
$ws.Range("J17").Value2 = 0.05133576516230916
$ws.Range("M17").Value2 = 67.71760166666667
$ws.Range("N17").Value2 = 203.152805
$ws.Range("O17").Value2 = 0.3089145827526977
$ws.Range("P17").Value2 = 0.3089145827526977
$ws.Range("Q17").Value2 = 3.528041901765556
$ws.Range("R17").Value2 = 31.75237711589
$ws.Range("S17").Value2 = 0.01585836647540521
$ws.Range("T17").Value2 = 0.01585836647540521
